$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 2380
$ws.Range("I12").Value = 4000
$ws.Range("J12").Value = 1300
$ws.Range("K12").Value = 4000
$ws.Range("L12").Value = 1300
$ws.Range("M12").Value = -3830
$ws.Range("N12").Value = -1640
# Row 33
$ws.Range("H33").Value = 306.5357
$ws.Range("I33").Value = 268.57693
$ws.Range("J33").Value = 800
$ws.Range("K33").Value = 268.57693
$ws.Range("L33").Value = 800
$ws.Range("M33").Value = -39.57693
$ws.Range("N33").Value = -1258
# Row 98
$ws.Range("H98").Value = 602.3
$ws.Range("I98").Value = 539.44446
$ws.Range("J98").Value = 1168
$ws.Range("K98").Value = 539.44446
$ws.Range("L98").Value = 1168
$ws.Range("M98").Value = 958.55554
$ws.Range("N98").Value = -4164
# Row 116
$ws.Range("H116").Value = 3842.2778
$ws.Range("I116").Value = 1764.1428
$ws.Range("K116").Value = 1764.1428
$ws.Range("M116").Value = 1677.8572
# Row 122
$ws.Range("H122").Value = 602.3
$ws.Range("I122").Value = 539.44446
$ws.Range("J122").Value = 1168
$ws.Range("K122").Value = 1618.33338
$ws.Range("L122").Value = 3504
$ws.Range("M122").Value = 831.66662
$ws.Range("N122").Value = -8404
# Row 129
$ws.Range("H129").Value = 792.931
$ws.Range("I129").Value = 647.125
$ws.Range("J129").Value = 848.4761999999999
$ws.Range("K129").Value = 1941.375
$ws.Range("L129").Value = 2545.4286
$ws.Range("M129").Value = 3058.625
$ws.Range("N129").Value = -12545.4286
# Row 132
$ws.Range("H132").Value = 2427.5938
$ws.Range("I132").Value = 2427.5938
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7282.7814
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4752.7814
$ws.Range("N132").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1928.1111
$ws.Range("I2").Value = 1238.6923
$ws.Range("J2").Value = 3720.6
$ws.Range("K2").Value = 1238.6923
$ws.Range("L2").Value = 3720.6
$ws.Range("M2").Value = -1125.6923
$ws.Range("N2").Value = -3946.6
# Row 45
$ws.Range("H45").Value = 2789.9285
$ws.Range("I45").Value = 2179.3333
$ws.Range("J45").Value = 3400.524
$ws.Range("K45").Value = 2179.3333
$ws.Range("L45").Value = 3400.524
$ws.Range("M45").Value = -1802.3333
$ws.Range("N45").Value = -4154.523999999999
# Row 74
$ws.Range("H74").Value = 2556.2856
$ws.Range("I74").Value = 2117.4
$ws.Range("J74").Value = 3653.5
$ws.Range("K74").Value = 2117.4
$ws.Range("L74").Value = 3653.5
$ws.Range("M74").Value = -1243.4
$ws.Range("N74").Value = -5401.5
# Row 77
$ws.Range("H77").Value = 2556.2856
$ws.Range("I77").Value = 2117.4
$ws.Range("J77").Value = 3653.5
$ws.Range("K77").Value = 10587
$ws.Range("L77").Value = 18267.5
$ws.Range("M77").Value = -6219
$ws.Range("N77").Value = -27003.5
# Row 102
$ws.Range("H102").Value = 2197.375
$ws.Range("I102").Value = 1217.5555
$ws.Range("J102").Value = 3457.1428
$ws.Range("K102").Value = 1217.5555
$ws.Range("L102").Value = 3457.1428
$ws.Range("M102").Value = 404.4445000000001
$ws.Range("N102").Value = -6701.1428
# Row 116
$ws.Range("H116").Value = 1928.1111
$ws.Range("I116").Value = 1238.6923
$ws.Range("J116").Value = 3720.6
$ws.Range("K116").Value = 1238.6923
$ws.Range("L116").Value = 3720.6
$ws.Range("M116").Value = 1055.3077
$ws.Range("N116").Value = -8308.6
# Row 132
$ws.Range("H132").Value = 37189.2
$ws.Range("I132").Value = 4541.3335
$ws.Range("J132").Value = 45351.168
$ws.Range("K132").Value = 13624.0005
$ws.Range("L132").Value = 136053.504
$ws.Range("M132").Value = -11094.0005
$ws.Range("N132").Value = -141113.504

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1928.1111
$ws.Range("I3").Value = 1238.6923
$ws.Range("J3").Value = 3720.6
$ws.Range("K3").Value = 1238.6923
$ws.Range("L3").Value = 3720.6
$ws.Range("M3").Value = -1124.6923
$ws.Range("N3").Value = -3948.6
# Row 99
$ws.Range("H99").Value = 2358.25
$ws.Range("I99").Value = 1929.8572
$ws.Range("J99").Value = 2958
$ws.Range("K99").Value = 1929.8572
$ws.Range("L99").Value = 2958
$ws.Range("M99").Value = -431.8571999999999
$ws.Range("N99").Value = -5954
# Row 134
$ws.Range("H134").Value = 76615.07000000001
$ws.Range("I134").Value = 76615.07000000001
$ws.Range("K134").Value = 229845.21
$ws.Range("M134").Value = -227310.21

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 13377.833
$ws.Range("I31").Value = 20125.4
$ws.Range("K31").Value = 20125.4
$ws.Range("M31").Value = -19830.4
# Row 34
$ws.Range("H34").Value = 13377.833
$ws.Range("I34").Value = 20125.4
$ws.Range("K34").Value = 20125.4
$ws.Range("M34").Value = -19923.4
# Row 107
$ws.Range("H107").Value = 555.125
$ws.Range("I107").Value = 416.0909
$ws.Range("J107").Value = 861
$ws.Range("K107").Value = 416.0909
$ws.Range("L107").Value = 861
$ws.Range("M107").Value = 1503.9091
$ws.Range("N107").Value = -4701
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
# Row 122
$ws.Range("H122").Value = 1002.3
$ws.Range("I122").Value = 1024.3
$ws.Range("J122").Value = 980.3
$ws.Range("K122").Value = 3072.9
$ws.Range("L122").Value = 2940.9
$ws.Range("M122").Value = -622.8999999999996
$ws.Range("N122").Value = -7840.9
# Row 132
$ws.Range("H132").Value = 20778.111
$ws.Range("I132").Value = 26053.5
$ws.Range("K132").Value = 78160.5
$ws.Range("M132").Value = -75630.5
# Row 134
$ws.Range("H134").Value = 1165.1333
$ws.Range("I134").Value = 997.7
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 2993.1
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -458.1000000000004
$ws.Range("N134").Value = -9570

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 775.27
$ws.Range("J131").Value = 780.5816
$ws.Range("L131").Value = 2341.7448
$ws.Range("N131").Value = -12421.7448

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 1657.75
$ws.Range("I107").Value = 256.66666
$ws.Range("J107").Value = 2498.4
$ws.Range("K107").Value = 256.66666
$ws.Range("L107").Value = 2498.4
$ws.Range("M107").Value = 1663.33334
$ws.Range("N107").Value = -6338.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3296.5
$ws.Range("I7").Value = 3042.5625
$ws.Range("J7").Value = 4312.25
$ws.Range("K7").Value = 3042.5625
$ws.Range("L7").Value = 4312.25
$ws.Range("M7").Value = -2930.5625
$ws.Range("N7").Value = -4536.25
# Row 22
$ws.Range("H22").Value = 5420.75
$ws.Range("I22").Value = 5401
$ws.Range("J22").Value = 5440.5
$ws.Range("K22").Value = 5401
$ws.Range("L22").Value = 5440.5
$ws.Range("M22").Value = -5106
$ws.Range("N22").Value = -6030.5
# Row 27
$ws.Range("H27").Value = 5420.75
$ws.Range("I27").Value = 5401
$ws.Range("J27").Value = 5440.5
$ws.Range("K27").Value = 5401
$ws.Range("L27").Value = 5440.5
$ws.Range("M27").Value = -5294
$ws.Range("N27").Value = -5654.5
# Row 46
$ws.Range("H46").Value = 1121.579
$ws.Range("I46").Value = 854.6667
$ws.Range("J46").Value = 2122.5
$ws.Range("K46").Value = 854.6667
$ws.Range("L46").Value = 2122.5
$ws.Range("M46").Value = -666.6667
$ws.Range("N46").Value = -2498.5
# Row 126
$ws.Range("H126").Value = 3296.5
$ws.Range("I126").Value = 3042.5625
$ws.Range("J126").Value = 4312.25
$ws.Range("K126").Value = 9127.6875
$ws.Range("L126").Value = 12936.75
$ws.Range("M126").Value = -6657.6875
$ws.Range("N126").Value = -17876.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 1608.7
$ws.Range("J107").Value = 1837.375
$ws.Range("L107").Value = 5512.125
$ws.Range("N107").Value = -9352.125
# Row 123
$ws.Range("H123").Value = 39429
$ws.Range("J123").Value = 39429
$ws.Range("L123").Value = 39429
$ws.Range("N123").Value = -49229
# Row 126
$ws.Range("H126").Value = 1073
$ws.Range("I126").Value = 681.3333
$ws.Range("J126").Value = 2080.1428
$ws.Range("K126").Value = 2043.9999
$ws.Range("L126").Value = 6240.428400000001
$ws.Range("M126").Value = 426.0001
$ws.Range("N126").Value = -11180.4284
# Row 132
$ws.Range("H132").Value = 2271.4517
$ws.Range("I132").Value = 2027.4762
$ws.Range("J132").Value = 2783.8
$ws.Range("K132").Value = 6082.4286
$ws.Range("L132").Value = 8351.400000000001
$ws.Range("M132").Value = -3552.4286
$ws.Range("N132").Value = -13411.4
# Row 136
$ws.Range("H136").Value = 1009256.5
$ws.Range("I136").Value = 1403374.1
$ws.Range("K136").Value = 4210122.300000001
$ws.Range("M136").Value = -4207572.300000001

Write-Output "Done applying changes"